$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 2.92
$ws.Range("G2").Value = 3.3
$ws.Range("H2").Value = 2.42
$ws.Range("I2").Value = 2.7
$ws.Range("J2").Value = 3.25
$ws.Range("O2").Value = 1.32
$ws.Range("P2").Value = 1.88
$ws.Range("R2").Value = 1.34
$ws.Range("T2").Value = 1.74
$ws.Range("V2").Value = 1.59
$ws.Range("W2").Value = 1.43
$ws.Range("Y2").Value = 970
$ws.Range("AE2").Value = 32
$ws.Range("AG2").Value = 970
$ws.Range("AI2").Value = 48
$ws.Range("AO2").Value = 26
$ws.Range("F3").Value = 2.56
$ws.Range("G3").Value = 2.58
$ws.Range("H3").Value = 2.66
$ws.Range("I3").Value = 2.86
$ws.Range("J3").Value = 3.8
$ws.Range("K3").Value = 4.2
$ws.Range("L3").Value = 1.33
$ws.Range("N3").Value = 2.74
$ws.Range("P3").Value = 1.81
$ws.Range("Q3").Value = 1.72
$ws.Range("R3").Value = 1.3
$ws.Range("S3").Value = 2.78
$ws.Range("V3").Value = 1.54
$ws.Range("W3").Value = 1.63
$ws.Range("X3").Value = 21
$ws.Range("Y3").Value = 14.5
$ws.Range("AC3").Value = 11
$ws.Range("AD3").Value = 14.5
$ws.Range("AH3").Value = 22
$ws.Range("F4").Value = 2.22
$ws.Range("G4").Value = 2.86
$ws.Range("H4").Value = 2.88
$ws.Range("I4").Value = 3.65
$ws.Range("J4").Value = 3.3
$ws.Range("K4").Value = 4.8
$ws.Range("P4").Value = 1.9
$ws.Range("Q4").Value = 1.67
$ws.Range("S4").Value = 2.6
$ws.Range("V4").Value = 1.38
$ws.Range("W4").Value = 1.53
$ws.Range("Z4").Value = 26
$ws.Range("AE4").Value = 38
$ws.Range("AF4").Value = 21
$ws.Range("AJ4").Value = 40
$ws.Range("AK4").Value = 30
$ws.Range("G5").Value = 8.800000000000001
$ws.Range("K5").Value = 5.4
$ws.Range("P5").Value = 2.34
$ws.Range("R5").Value = 1.52
$ws.Range("S5").Value = 2.86
$ws.Range("T5").Value = 1.94
$ws.Range("AJ5").Value = 290
$ws.Range("H6").Value = 1.77
$ws.Range("N6").Value = 5.5
$ws.Range("P6").Value = 2.56
$ws.Range("Q6").Value = 1.54
$ws.Range("R6").Value = 1.63
$ws.Range("S6").Value = 2.28
$ws.Range("T6").Value = 1.58
$ws.Range("U6").Value = 2.46
$ws.Range("W6").Value = 1.26
$ws.Range("X6").Value = 28
$ws.Range("Y6").Value = 13.5
$ws.Range("AC6").Value = 11.5
$ws.Range("AE6").Value = 19.5
$ws.Range("AF6").Value = 46
$ws.Range("AG6").Value = 22
$ws.Range("AH6").Value = 18.5
$ws.Range("AI6").Value = 28
$ws.Range("AO6").Value = 9
$ws.Range("I7").Value = 1.89
$ws.Range("M7").Value = 1.05
$ws.Range("Q7").Value = 1.76
$ws.Range("R7").Value = 1.49
$ws.Range("AB7").Value = 18.5
$ws.Range("AG7").Value = 17
$ws.Range("AI7").Value = 29
$ws.Range("F8").Value = 1.18
$ws.Range("G8").Value = 1.69
$ws.Range("J8").Value = 1.09
$ws.Range("P8").Value = 1.09
$ws.Range("J9").Value = 3.45
$ws.Range("K9").Value = 3.5
$ws.Range("L9").Value = 1.44
$ws.Range("N9").Value = 3.75
$ws.Range("AA9").Value = 60
$ws.Range("AF9").Value = 14
$ws.Range("AH9").Value = 17.5
$ws.Range("AI9").Value = 50
$ws.Range("AK9").Value = 25
$ws.Range("AM9").Value = 95
$ws.Range("AN9").Value = 20
$ws.Range("L10").Value = 1.34
$ws.Range("O10").Value = 1.24
$ws.Range("P10").Value = 2.26
$ws.Range("R10").Value = 1.5
$ws.Range("S10").Value = 2.96
$ws.Range("Y10").Value = 8.199999999999999
$ws.Range("Z10").Value = 7
$ws.Range("F11").Value = 2.7
$ws.Range("G11").Value = 2.72
$ws.Range("H11").Value = 2.88
$ws.Range("L11").Value = 1.41
$ws.Range("N11").Value = 4
$ws.Range("P11").Value = 1.99
$ws.Range("Q11").Value = 1.99
$ws.Range("W11").Value = 1.58
$ws.Range("AI11").Value = 42
$ws.Range("H12").Value = 1.96
$ws.Range("I12").Value = 1.98
$ws.Range("L12").Value = 1.33
$ws.Range("N12").Value = 5.1
$ws.Range("X12").Value = 19
$ws.Range("Y12").Value = 11.5
$ws.Range("AJ12").Value = 75
$ws.Range("G13").Value = 4.9
$ws.Range("J13").Value = 4.2
$ws.Range("L13").Value = 1.29
$ws.Range("N13").Value = 5.9
$ws.Range("O13").Value = 1.19
$ws.Range("P13").Value = 2.62